$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text (matches source inlineStr
# cells), bypassing Excel automatic number/date inference, and without leaving
# a residual cell style behind.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.792.03"
Set-TextValue "E2" "  +0.70%  "
Set-TextValue "D3" "2.318.38"
Set-TextValue "E3" "  +4.39%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "97.24"
Set-TextValue "E5" "  +5.23%  "
Set-TextValue "D6" "270.70"
Set-TextValue "E6" "  +0.38%  "
Set-TextValue "E7" "  +0.56%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.623"
Set-TextValue "E9" "  +0.90%  "
Set-TextValue "D10" "45.59"
Set-TextValue "E10" "  -0.65%  "
Set-TextValue "D11" "0.0948"
Set-TextValue "E11" "  +1.42%  "
Set-TextValue "D12" "8.04"
Set-TextValue "E12" "  -1.79%  "
Set-TextValue "E13" "  +0.52%  "
Set-TextValue "D14" "2.660.01"
Set-TextValue "E14" "  +4.03%  "
Set-TextValue "E15" "  +3.16%  "
Set-TextValue "D16" "0.868"
Set-TextValue "E16" "  +8.49%  "
Set-TextValue "D17" "2.321.87"
Set-TextValue "E17" "  +3.67%  "
Set-TextValue "D18" "43.728.51"
Set-TextValue "E18" "  +0.56%  "
Set-TextValue "E19" "  +5.38%  "
Set-TextValue "D20" "6.39"
Set-TextValue "E20" "  +6.46%  "
Set-TextValue "D21" "72.72"
Set-TextValue "E21" "  +3.34%  "
Set-TextValue "D22" "239.57"
Set-TextValue "E22" "  +3.11%  "
Set-TextValue "D23" "2.27"
Set-TextValue "E23" "  -3.35%  "
Set-TextValue "D24" "9.36"
Set-TextValue "E24" "  +3.77%  "
Set-TextValue "E25" "  -0.08%  "
Set-TextValue "E26" "  +1.14%  "
Set-TextValue "D27" "11.34"
Set-TextValue "E27" "  +0.47%  "
Set-TextValue "E28" "  -2.00%  "
Set-TextValue "D29" "2.28"
Set-TextValue "E29" "  +0.84%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "22.41"
Set-TextValue "E30" "  +7.72%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D31" "37.97"
Set-TextValue "E31" "  -8.65%  "
Set-TextValue "D32" "175.06"
Set-TextValue "E32" "  +1.40%  "
Set-TextValue "D33" "0.0899"
Set-TextValue "E33" "  -2.23%  "
Set-TextValue "D34" "5.47"
Set-TextValue "E34" "  +0.38%  "
Set-TextValue "E35" "  +3.18%  "
Set-TextValue "D36" "0.0361"
Set-TextValue "E36" "  +3.13%  "
Set-TextValue "E37" "  -2.95%  "
Set-TextValue "D38" "4.38"
Set-TextValue "E38" "  +1.89%  "
Set-TextValue "E39" "  -6.61%  "
Set-TextValue "E40" "  +11.52%  "
Set-TextValue "D41" "2.35"
Set-TextValue "E41" "  +9.20%  "
Set-TextValue "D42" "1.37"
Set-TextValue "E42" "  +19.20%  "
Set-TextValue "E43" "  -3.08%  "
Set-TextValue "D44" "9.17"
Set-TextValue "E44" "  +9.71%  "
Set-TextValue "D45" "62.09"
Set-TextValue "E45" "  -1.72%  "
Set-TextValue "D46" "5.35"
Set-TextValue "E46" "  +0.77%  "
Set-TextValue "E47" "  +4.38%  "
Set-TextValue "D48" "100.34"
Set-TextValue "E48" "  +0.22%  "
Set-TextValue "E49" "  +0.65%  "
Set-TextValue "D50" "0.191"
Set-TextValue "E50" "  +18.15%  "
Set-TextValue "D51" "2.546.56"
Set-TextValue "E51" "  +4.14%  "
